# Auto-generated script applying scheduled-runner market data updates
# to the Kujata_Profits workbook (per-sheet "Leve" profit tables).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(33, 8).Value = 124.666664   # H33: 111.63158 -> 124.666664
$ws.Cells.Item(33, 9).Value = 103.76471   # I33: 91.166664 -> 103.76471
$ws.Cells.Item(33, 11).Value = 103.76471   # K33: 91.166664 -> 103.76471
$ws.Cells.Item(33, 13).Value = 125.23529   # M33: 137.833336 -> 125.23529
$ws.Cells.Item(40, 8).Value = 1788.9   # H40: 1691.5 -> 1788.9
$ws.Cells.Item(40, 9).Value = 1900   # I40: 1533 -> 1900
$ws.Cells.Item(40, 10).Value = 1761.125   # J40: 1744.3334 -> 1761.125
$ws.Cells.Item(40, 11).Value = 1900   # K40: 1533 -> 1900
$ws.Cells.Item(40, 12).Value = 1761.125   # L40: 1744.3334 -> 1761.125
$ws.Cells.Item(40, 13).Value = -1725   # M40: -1358 -> -1725
$ws.Cells.Item(40, 14).Value = -2111.125   # N40: -2094.3334 -> -2111.125
$ws.Cells.Item(93, 8).Value = 20601   # H93: 30000 -> 20601
$ws.Cells.Item(93, 10).Value = 20601   # J93: 30000 -> 20601
$ws.Cells.Item(93, 12).Value = 20601   # L93: 30000 -> 20601
$ws.Cells.Item(93, 14).Value = -25593   # N93: -34992 -> -25593
$ws.Cells.Item(137, 8).Value = 1430.2142   # H137: 1503 -> 1430.2142
$ws.Cells.Item(137, 9).Value = 1368.3334   # I137: 1493 -> 1368.3334
$ws.Cells.Item(137, 10).Value = 1801.5   # J137: 1603 -> 1801.5
$ws.Cells.Item(137, 11).Value = 4105.0002   # K137: 4479 -> 4105.0002
$ws.Cells.Item(137, 12).Value = 5404.5   # L137: 4809 -> 5404.5
$ws.Cells.Item(137, 13).Value = -1555.0002   # M137: -1929 -> -1555.0002
$ws.Cells.Item(137, 14).Value = -10504.5   # N137: -9909 -> -10504.5

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 3981.2744   # H32: 3117.5 -> 3981.2744
$ws.Cells.Item(32, 9).Value = 3587.3914   # I32: 2725.15 -> 3587.3914
$ws.Cells.Item(32, 10).Value = 7605   # J32: 9002.75 -> 7605
$ws.Cells.Item(32, 11).Value = 3587.3914   # K32: 2725.15 -> 3587.3914
$ws.Cells.Item(32, 12).Value = 7605   # L32: 9002.75 -> 7605
$ws.Cells.Item(32, 13).Value = -3300.3914   # M32: -2438.15 -> -3300.3914
$ws.Cells.Item(32, 14).Value = -8179   # N32: -9576.75 -> -8179
$ws.Cells.Item(96, 8).Value = 17399.6   # H96: 17599.8 -> 17399.6
$ws.Cells.Item(96, 10).Value = 17399.6   # J96: 17599.8 -> 17399.6
$ws.Cells.Item(96, 12).Value = 17399.6   # L96: 17599.8 -> 17399.6
$ws.Cells.Item(96, 14).Value = -22891.6   # N96: -23091.8 -> -22891.6
$ws.Cells.Item(102, 8).Value = 33335538   # H102: 23811384 -> 33335538
$ws.Cells.Item(102, 9).Value = 33335538   # I102: 27779782 -> 33335538
$ws.Cells.Item(102, 10).Value = 0   # J102: 1000 -> 0
$ws.Cells.Item(102, 11).Value = 33335538   # K102: 27779782 -> 33335538
$ws.Cells.Item(102, 12).Value = 0   # L102: 1000 -> 0
$ws.Cells.Item(102, 13).Value = -33333916   # M102: -27778160 -> -33333916
$ws.Cells.Item(102, 14).Value = $null   # N102: remove (was -4244)
$ws.Cells.Item(119, 8).Value = 31698   # H119: 0 -> 31698
$ws.Cells.Item(119, 10).Value = 31698   # J119: 0 -> 31698
$ws.Cells.Item(119, 12).Value = 31698   # L119: 0 -> 31698
$ws.Cells.Item(119, 14).Value = -41374   # N119: None -> -41374
$ws.Cells.Item(122, 8).Value = 1155.3158   # H122: 1398.9 -> 1155.3158
$ws.Cells.Item(122, 9).Value = 1114.8235   # I122: 1425.1428 -> 1114.8235
$ws.Cells.Item(122, 10).Value = 1499.5   # J122: 1337.6666 -> 1499.5
$ws.Cells.Item(122, 11).Value = 3344.4705   # K122: 4275.428400000001 -> 3344.4705
$ws.Cells.Item(122, 12).Value = 4498.5   # L122: 4012.9998 -> 4498.5
$ws.Cells.Item(122, 13).Value = -894.4704999999999   # M122: -1825.428400000001 -> -894.4704999999999
$ws.Cells.Item(122, 14).Value = -9398.5   # N122: -8912.9998 -> -9398.5
$ws.Cells.Item(132, 8).Value = 1928.5   # H132: 1977.4839 -> 1928.5
$ws.Cells.Item(132, 9).Value = 1530.6666   # I132: 1536.5 -> 1530.6666
$ws.Cells.Item(132, 10).Value = 3122   # J132: 3489.4285 -> 3122
$ws.Cells.Item(132, 11).Value = 4591.9998   # K132: 4609.5 -> 4591.9998
$ws.Cells.Item(132, 12).Value = 9366   # L132: 10468.2855 -> 9366
$ws.Cells.Item(132, 13).Value = -2061.9998   # M132: -2079.5 -> -2061.9998
$ws.Cells.Item(132, 14).Value = -14426   # N132: -15528.2855 -> -14426
$ws.Cells.Item(133, 8).Value = 28889.445   # H133: 28296.666 -> 28889.445
$ws.Cells.Item(133, 9).Value = 0   # I133: 28000 -> 0
$ws.Cells.Item(133, 10).Value = 28889.445   # J133: 28309.564 -> 28889.445
$ws.Cells.Item(133, 11).Value = 0   # K133: 28000 -> 0
$ws.Cells.Item(133, 12).Value = 28889.445   # L133: 28309.564 -> 28889.445
$ws.Cells.Item(133, 13).Value = $null   # M133: remove (was -25470)
$ws.Cells.Item(133, 14).Value = -33949.445   # N133: -33369.564 -> -33949.445

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(46, 8).Value = 0   # H46: 3800 -> 0
$ws.Cells.Item(46, 9).Value = 0   # I46: 3800 -> 0
$ws.Cells.Item(46, 11).Value = 0   # K46: 3800 -> 0
$ws.Cells.Item(46, 13).Value = $null   # M46: remove (was -3502)
$ws.Cells.Item(103, 8).Value = 0   # H103: 10657 -> 0
$ws.Cells.Item(103, 10).Value = 0   # J103: 10657 -> 0
$ws.Cells.Item(103, 12).Value = 0   # L103: 10657 -> 0
$ws.Cells.Item(103, 14).Value = $null   # N103: remove (was -13001)

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 8).Value = 1172.125   # H6: 1375.25 -> 1172.125
$ws.Cells.Item(6, 10).Value = 1075.2   # J6: 1500 -> 1075.2
$ws.Cells.Item(6, 12).Value = 1075.2   # L6: 1500 -> 1075.2
$ws.Cells.Item(6, 14).Value = -1301.2   # N6: -1726 -> -1301.2
$ws.Cells.Item(7, 8).Value = 129.2   # H7: 107.833336 -> 129.2
$ws.Cells.Item(7, 9).Value = 77.666664   # I7: 58.5 -> 77.666664
$ws.Cells.Item(7, 11).Value = 77.666664   # K7: 58.5 -> 77.666664
$ws.Cells.Item(7, 13).Value = 35.333336   # M7: 54.5 -> 35.333336
$ws.Cells.Item(62, 8).Value = 7410251.5   # H62: 8002852 -> 7410251.5
$ws.Cells.Item(62, 9).Value = 2971.76   # I62: 2970.5833 -> 2971.76
$ws.Cells.Item(62, 10).Value = 100001250   # J62: 200000000 -> 100001250
$ws.Cells.Item(62, 11).Value = 2971.76   # K62: 2970.5833 -> 2971.76
$ws.Cells.Item(62, 12).Value = 100001250   # L62: 200000000 -> 100001250
$ws.Cells.Item(62, 13).Value = -2347.76   # M62: -2346.5833 -> -2347.76
$ws.Cells.Item(62, 14).Value = -100002498   # N62: -200001248 -> -100002498
$ws.Cells.Item(65, 8).Value = 7410251.5   # H65: 8002852 -> 7410251.5
$ws.Cells.Item(65, 9).Value = 2971.76   # I65: 2970.5833 -> 2971.76
$ws.Cells.Item(65, 10).Value = 100001250   # J65: 200000000 -> 100001250
$ws.Cells.Item(65, 11).Value = 14858.8   # K65: 14852.9165 -> 14858.8
$ws.Cells.Item(65, 12).Value = 500006250   # L65: 1000000000 -> 500006250
$ws.Cells.Item(65, 13).Value = -11738.8   # M65: -11732.9165 -> -11738.8
$ws.Cells.Item(65, 14).Value = -500012490   # N65: -1000006240 -> -500012490
$ws.Cells.Item(95, 8).Value = 5184.6665   # H95: 7839.5 -> 5184.6665
$ws.Cells.Item(95, 10).Value = 5184.6665   # J95: 7839.5 -> 5184.6665
$ws.Cells.Item(95, 12).Value = 5184.6665   # L95: 7839.5 -> 5184.6665
$ws.Cells.Item(95, 14).Value = -10676.6665   # N95: -13331.5 -> -10676.6665
$ws.Cells.Item(132, 8).Value = 6353   # H132: 5288.533 -> 6353
$ws.Cells.Item(132, 9).Value = 7193.278   # I132: 5652.625 -> 7193.278
$ws.Cells.Item(132, 11).Value = 21579.834   # K132: 16957.875 -> 21579.834
$ws.Cells.Item(132, 13).Value = -19049.834   # M132: -14427.875 -> -19049.834
$ws.Cells.Item(134, 8).Value = 1739.2812   # H134: 1665.7428 -> 1739.2812
$ws.Cells.Item(134, 9).Value = 1746.625   # I134: 1687.0385 -> 1746.625
$ws.Cells.Item(134, 10).Value = 1717.25   # J134: 1604.2222 -> 1717.25
$ws.Cells.Item(134, 11).Value = 5239.875   # K134: 5061.1155 -> 5239.875
$ws.Cells.Item(134, 12).Value = 5151.75   # L134: 4812.6666 -> 5151.75
$ws.Cells.Item(134, 13).Value = -2704.875   # M134: -2526.1155 -> -2704.875
$ws.Cells.Item(134, 14).Value = -10221.75   # N134: -9882.6666 -> -10221.75

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(8, 8).Value = 62.22222   # H8: 142929.28 -> 62.22222
$ws.Cells.Item(8, 9).Value = 62.22222   # I8: 142929.28 -> 62.22222
$ws.Cells.Item(8, 11).Value = 186.66666   # K8: 428787.84 -> 186.66666
$ws.Cells.Item(8, 13).Value = -47.66666000000001   # M8: -428648.84 -> -47.66666000000001
$ws.Cells.Item(39, 8).Value = 1699.4546   # H39: 1653.3846 -> 1699.4546
$ws.Cells.Item(39, 10).Value = 1699.4546   # J39: 1653.3846 -> 1699.4546
$ws.Cells.Item(39, 12).Value = 5098.3638   # L39: 4960.1538 -> 5098.3638
$ws.Cells.Item(39, 14).Value = -5686.3638   # N39: -5548.1538 -> -5686.3638
$ws.Cells.Item(114, 8).Value = 609.6667   # H114: 578.7 -> 609.6667
$ws.Cells.Item(114, 10).Value = 720.5   # J114: 580.3333 -> 720.5
$ws.Cells.Item(114, 12).Value = 2161.5   # L114: 1740.9999 -> 2161.5
$ws.Cells.Item(114, 14).Value = -8669.5   # N114: -8248.999900000001 -> -8669.5
$ws.Cells.Item(132, 8).Value = 846.6667   # H132: 911.25 -> 846.6667
$ws.Cells.Item(132, 9).Value = 490   # I132: 500 -> 490
$ws.Cells.Item(132, 10).Value = 1025   # J132: 970 -> 1025
$ws.Cells.Item(132, 11).Value = 4410   # K132: 4500 -> 4410
$ws.Cells.Item(132, 12).Value = 9225   # L132: 8730 -> 9225
$ws.Cells.Item(132, 13).Value = -1880   # M132: -1970 -> -1880
$ws.Cells.Item(132, 14).Value = -14285   # N132: -13790 -> -14285

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(64, 8).Value = 16500   # H64: 0 -> 16500
$ws.Cells.Item(64, 10).Value = 16500   # J64: 0 -> 16500
$ws.Cells.Item(64, 12).Value = 16500   # L64: 0 -> 16500
$ws.Cells.Item(64, 14).Value = -16996   # N64: None -> -16996
$ws.Cells.Item(67, 8).Value = 16500   # H67: 0 -> 16500
$ws.Cells.Item(67, 10).Value = 16500   # J67: 0 -> 16500
$ws.Cells.Item(67, 12).Value = 16500   # L67: 0 -> 16500
$ws.Cells.Item(67, 14).Value = -18216   # N67: None -> -18216
$ws.Cells.Item(70, 8).Value = 56256250   # H70: 56255252 -> 56256250
$ws.Cells.Item(70, 9).Value = 62505624   # I70: 50004900 -> 62505624
$ws.Cells.Item(70, 10).Value = 50006876   # J70: 66672500 -> 50006876
$ws.Cells.Item(70, 11).Value = 62505624   # K70: 50004900 -> 62505624
$ws.Cells.Item(70, 12).Value = 50006876   # L70: 66672500 -> 50006876
$ws.Cells.Item(70, 13).Value = -62505354   # M70: -50004630 -> -62505354
$ws.Cells.Item(70, 14).Value = -50007416   # N70: -66673040 -> -50007416
$ws.Cells.Item(73, 8).Value = 56256250   # H73: 56255252 -> 56256250
$ws.Cells.Item(73, 9).Value = 62505624   # I73: 50004900 -> 62505624
$ws.Cells.Item(73, 10).Value = 50006876   # J73: 66672500 -> 50006876
$ws.Cells.Item(73, 11).Value = 62505624   # K73: 50004900 -> 62505624
$ws.Cells.Item(73, 12).Value = 50006876   # L73: 66672500 -> 50006876
$ws.Cells.Item(73, 13).Value = -62504688   # M73: -50003964 -> -62504688
$ws.Cells.Item(73, 14).Value = -50008748   # N73: -66674372 -> -50008748
$ws.Cells.Item(95, 8).Value = 18372.5   # H95: 18622.5 -> 18372.5
$ws.Cells.Item(95, 10).Value = 18372.5   # J95: 18622.5 -> 18372.5
$ws.Cells.Item(95, 12).Value = 18372.5   # L95: 18622.5 -> 18372.5
$ws.Cells.Item(95, 14).Value = -23864.5   # N95: -24114.5 -> -23864.5
$ws.Cells.Item(126, 8).Value = 2014.05   # H126: 2140.0588 -> 2014.05
$ws.Cells.Item(126, 9).Value = 1961.5454   # I126: 2107.7 -> 1961.5454
$ws.Cells.Item(126, 10).Value = 2078.2222   # J126: 2186.2856 -> 2078.2222
$ws.Cells.Item(126, 11).Value = 5884.6362   # K126: 6323.099999999999 -> 5884.6362
$ws.Cells.Item(126, 12).Value = 6234.6666   # L126: 6558.8568 -> 6234.6666
$ws.Cells.Item(126, 13).Value = -3414.6362   # M126: -3853.099999999999 -> -3414.6362
$ws.Cells.Item(126, 14).Value = -11174.6666   # N126: -11498.8568 -> -11174.6666
$ws.Cells.Item(132, 8).Value = 3084.0667   # H132: 4699.5 -> 3084.0667
$ws.Cells.Item(132, 9).Value = 2751.182   # I132: 6099.5 -> 2751.182
$ws.Cells.Item(132, 11).Value = 8253.545999999998   # K132: 18298.5 -> 8253.545999999998
$ws.Cells.Item(132, 13).Value = -5723.545999999998   # M132: -15768.5 -> -5723.545999999998

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 1311.25   # H7: 1276.6666 -> 1311.25
$ws.Cells.Item(7, 9).Value = 1140.7142   # I7: 1123.125 -> 1140.7142
$ws.Cells.Item(7, 11).Value = 1140.7142   # K7: 1123.125 -> 1140.7142
$ws.Cells.Item(7, 13).Value = -1028.7142   # M7: -1011.125 -> -1028.7142
$ws.Cells.Item(94, 8).Value = 13199.4   # H94: 12499.75 -> 13199.4
$ws.Cells.Item(94, 10).Value = 13199.4   # J94: 12499.75 -> 13199.4
$ws.Cells.Item(94, 12).Value = 13199.4   # L94: 12499.75 -> 13199.4
$ws.Cells.Item(94, 14).Value = -14551.4   # N94: -13851.75 -> -14551.4
$ws.Cells.Item(100, 8).Value = 0   # H100: 2441.5 -> 0
$ws.Cells.Item(100, 9).Value = 0   # I100: 2133 -> 0
$ws.Cells.Item(100, 10).Value = 0   # J100: 2750 -> 0
$ws.Cells.Item(100, 11).Value = 0   # K100: 2133 -> 0
$ws.Cells.Item(100, 12).Value = 0   # L100: 2750 -> 0
$ws.Cells.Item(100, 13).Value = $null   # M100: remove (was -1592)
$ws.Cells.Item(100, 14).Value = $null   # N100: remove (was -3832)
$ws.Cells.Item(110, 8).Value = 29257.4   # H110: 29660.75 -> 29257.4
$ws.Cells.Item(110, 10).Value = 29071.75   # J110: 29547.666 -> 29071.75
$ws.Cells.Item(110, 12).Value = 29071.75   # L110: 29547.666 -> 29071.75
$ws.Cells.Item(110, 14).Value = -37251.75   # N110: -37727.666 -> -37251.75
$ws.Cells.Item(122, 8).Value = 22729936   # H122: 35717412 -> 22729936
$ws.Cells.Item(122, 9).Value = 31252550   # I122: 62503252 -> 31252550
$ws.Cells.Item(122, 11).Value = 93757650   # K122: 187509756 -> 93757650
$ws.Cells.Item(122, 13).Value = -93755200   # M122: -187507306 -> -93755200
$ws.Cells.Item(126, 8).Value = 1311.25   # H126: 1276.6666 -> 1311.25
$ws.Cells.Item(126, 9).Value = 1140.7142   # I126: 1123.125 -> 1140.7142
$ws.Cells.Item(126, 11).Value = 3422.1426   # K126: 3369.375 -> 3422.1426
$ws.Cells.Item(126, 13).Value = -952.1425999999997   # M126: -899.375 -> -952.1425999999997
$ws.Cells.Item(133, 8).Value = 34632.145   # H133: 35824.625 -> 34632.145
$ws.Cells.Item(133, 10).Value = 34632.145   # J133: 35824.625 -> 34632.145
$ws.Cells.Item(133, 12).Value = 34632.145   # L133: 35824.625 -> 34632.145
$ws.Cells.Item(133, 14).Value = -39692.145   # N133: -40884.625 -> -39692.145
$ws.Cells.Item(136, 8).Value = 5592.7036   # H136: 5761.654 -> 5592.7036
$ws.Cells.Item(136, 9).Value = 8606.071   # I136: 9175.77 -> 8606.071
$ws.Cells.Item(136, 11).Value = 25818.213   # K136: 27527.31 -> 25818.213
$ws.Cells.Item(136, 13).Value = -23268.213   # M136: -24977.31 -> -23268.213
$ws.Cells.Item(140, 8).Value = 38813.332   # H140: 58400 -> 38813.332
$ws.Cells.Item(140, 10).Value = 38813.332   # J140: 58400 -> 38813.332
$ws.Cells.Item(140, 12).Value = 38813.332   # L140: 58400 -> 38813.332
$ws.Cells.Item(140, 14).Value = -49173.332   # N140: -68760 -> -49173.332

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(100, 8).Value = 1250   # H100: 434.8 -> 1250
$ws.Cells.Item(100, 9).Value = 1250   # I100: 393.5 -> 1250
$ws.Cells.Item(100, 10).Value = 0   # J100: 600 -> 0
$ws.Cells.Item(100, 11).Value = 2500   # K100: 787 -> 2500
$ws.Cells.Item(100, 12).Value = 0   # L100: 1200 -> 0
$ws.Cells.Item(100, 13).Value = -1959   # M100: -246 -> -1959
$ws.Cells.Item(100, 14).Value = $null   # N100: remove (was -2282)
$ws.Cells.Item(119, 8).Value = 50698   # H119: 0 -> 50698
$ws.Cells.Item(119, 10).Value = 50698   # J119: 0 -> 50698
$ws.Cells.Item(119, 12).Value = 50698   # L119: 0 -> 50698
$ws.Cells.Item(119, 14).Value = -60374   # N119: None -> -60374
$ws.Cells.Item(132, 8).Value = 2741.516   # H132: 2666.1562 -> 2741.516
$ws.Cells.Item(132, 9).Value = 2330.2173   # I132: 2246.875 -> 2330.2173
$ws.Cells.Item(132, 11).Value = 6990.651899999999   # K132: 6740.625 -> 6990.651899999999
$ws.Cells.Item(132, 13).Value = -4460.651899999999   # M132: -4210.625 -> -4460.651899999999
